$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2021-12-31"

# Row 10 (August) - 2020 (Q,R,S) and 2021 (T,U,V) updates
$ws.Range("Q10").Value = 8
$ws.Range("R10").Value = 155
$ws.Range("S10").Value = 0.0491
$ws.Range("T10").Value = 10
$ws.Range("U10").Value = 150
$ws.Range("V10").Value = 0.0625

# Row 11 (September) - 2021 (T,U,V) updates
$ws.Range("T11").Value = 8
$ws.Range("U11").Value = 169
$ws.Range("V11").Value = 0.0452

# Row 13 (November) - 2021 (T,U,V) updates
$ws.Range("T13").Value = 7
$ws.Range("U13").Value = 194
$ws.Range("V13").Value = 0.0348

# Row 14 (December) - label change and all years updates
$ws.Range("A14").Value = "December (through 12-31)"

$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 44
$ws.Range("D14").Value = 0.102
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 91
$ws.Range("G14").Value = 0.09
$ws.Range("H14").Value = 13
$ws.Range("I14").Value = 103
$ws.Range("J14").Value = 0.1121
$ws.Range("K14").Value = 5
$ws.Range("L14").Value = 77
$ws.Range("M14").Value = 0.061
$ws.Range("N14").Value = 5
$ws.Range("O14").Value = 64
$ws.Range("P14").Value = 0.0725
$ws.Range("Q14").Value = 9
$ws.Range("R14").Value = 140
$ws.Range("S14").Value = 0.0604
$ws.Range("T14").Value = 2
$ws.Range("U14").Value = 202
$ws.Range("V14").Value = 0.0098

# Row 15 (Total) - all years updates
$ws.Range("B15").Value = 38
$ws.Range("C15").Value = 302
$ws.Range("D15").Value = 0.1118
$ws.Range("E15").Value = 68
$ws.Range("F15").Value = 595
$ws.Range("G15").Value = 0.1026
$ws.Range("H15").Value = 76
$ws.Range("I15").Value = 861
$ws.Range("J15").Value = 0.0811
$ws.Range("K15").Value = 79
$ws.Range("L15").Value = 685
$ws.Range("M15").Value = 0.1034
$ws.Range("N15").Value = 59
$ws.Range("O15").Value = 544
$ws.Range("P15").Value = 0.0978
$ws.Range("Q15").Value = 74
$ws.Range("R15").Value = 1339
$ws.Range("S15").Value = 0.0524
$ws.Range("T15").Value = 106
$ws.Range("U15").Value = 1742
$ws.Range("V15").Value = 0.0574
